# Refresh the "cryptos" price/volume snapshot (Price = column D, Volume(1h) = column E).
# Note: some Price values (e.g. "6.13") read as plain decimals, so a leading
# apostrophe is used to force them to stay text instead of being auto-converted
# to a number by Excel, matching the other Price cells which are already text
# (some, like "64.382.76", have multiple dots and are never auto-numeric).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.382.76'
$ws.Range("E2").Value = '  +1.41%  '
$ws.Range("D3").Value = '3.097.93'
$ws.Range("E3").Value = '  +1.02%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''559.69'
$ws.Range("E5").Value = '  +1.72%  '
$ws.Range("D6").Value = '''144.59'
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.094.67'
$ws.Range("E8").Value = '  +1.08%  '
$ws.Range("E9").Value = '  +0.77%  '
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("D11").Value = '''6.13'
$ws.Range("E11").Value = '  -6.49%  '
$ws.Range("E12").Value = '  +3.33%  '
$ws.Range("D13").Value = '''0.0000228'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("D14").Value = '''35.22'
$ws.Range("E14").Value = '  +0.68%  '
$ws.Range("D15").Value = '3.593.04'
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("D16").Value = '64.403.69'
$ws.Range("E16").Value = '  +1.41%  '
$ws.Range("D17").Value = '3.090.52'
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("D18").Value = '''0.111'
$ws.Range("E18").Value = '  +1.31%  '
$ws.Range("D20").Value = '''485.61'
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("D21").Value = '''14.02'
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").Value = '''7.57'
$ws.Range("E23").Value = '  +4.10%  '
$ws.Range("D24").Value = '''14.25'
$ws.Range("E24").Value = '  +12.49%  '
$ws.Range("D25").Value = '''81.32'
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '''2.81'
$ws.Range("E27").Value = '  +1.25%  '
$ws.Range("D28").Value = '''8.03'
$ws.Range("E28").Value = '  +1.52%  '
$ws.Range("E29").Value = '  +2.47%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("D32").Value = '''1.14'
$ws.Range("E32").Value = '  -1.53%  '
$ws.Range("E33").Value = '  +0.61%  '
$ws.Range("D34").Value = '''5.61'
$ws.Range("E34").Value = '  -1.50%  '
$ws.Range("D35").Value = '''6.23'
$ws.Range("E35").Value = '  +3.92%  '
$ws.Range("D36").Value = '''55.71'
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").Value = '''3.00'
$ws.Range("E37").Value = '  +17.10%  '
$ws.Range("D38").Value = '''452.79'
$ws.Range("E38").Value = '  -2.71%  '
$ws.Range("D39").Value = '''0.0410'
$ws.Range("E39").Value = '  +2.93%  '
$ws.Range("D40").Value = '''0.0817'
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("D41").Value = '2.965.03'
$ws.Range("E41").Value = '  -3.29%  '
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("E43").Value = '  -5.23%  '
$ws.Range("D44").Value = '''28.26'
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("E45").Value = '  +1.79%  '
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("E47").Value = '  +3.82%  '
$ws.Range("E48").Value = '  +1.73%  '
$ws.Range("D49").Value = '''118.77'
$ws.Range("E49").Value = '  +1.45%  '
$ws.Range("D50").Value = '0.0₃0517'
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("E51").Value = '  +0.37%  '
